# The underlying OOXML diff for this revision (see commit message: "Fixed
# POI packaging and upgraded to POI 3.15") is purely a re-serialization
# artifact: every changed line is the exact same element with the exact
# same set of attribute names/values, just re-ordered (mostly
# alphabetized) by the XML writer of the newer Apache POI version that
# regenerated this test fixture. There is no textual, formatting, or
# structural change to the document's content anywhere in document.xml
# or styles.xml - paragraphs, runs, tab stops, page size/margins, run
# fonts, language, latent style table, and style definitions all keep
# identical values before and after, only attribute order differs.
#
# Word's object model (and this COM-interop runtime) does not expose any
# way to control the attribute-serialization order used when it writes
# OOXML - that is purely an artifact of whichever library produced the
# file, not something reachable via Document/Range/Paragraph/Find
# COM calls. Since there is no semantic change for Word to make, the
# faithful reproduction of this revision is simply to leave the
# document's content untouched.
$d = $word.ActiveDocument
